$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "[name=`"Spokesman Czarny`"]  As a result, the girl of today was pushed hastily into the position, our 'youngest Platinum.'`n"

$ws.Range("C42").Value = "[name=`"Bald Marcin`"]  Against 'Left-hand' Tytus Topola.`n"

$ws.Range("C67").Value = "[name=`"Greatmouth Mob`"]  And facing her—! The concealed blade, the height of status, named for the old knight of legend, it’s the headliner of the Blade Helmet Knightclub, it’s 'Left-hand' Tytus Topola! `n"

$ws.Range("C72").Value = "[name=`"Greatmouth Mob`"]  After defeating 'Plastic' Szewczyk, she’s been marching victorious! Big or small, no matter the event, Maria Nearl brings her A-game!`n"

$ws.Range("C113").Value = "In the shadow of 'Left-hand' Tytus, Maria raises her sword once more.`n"
